# "Generate Report for Handoff"
# Updates the localization-status report: the two target-language sheets
# (zh-cn, de-de) now report Priority "ht" for the rows that are Ready for
# handoff, and the handoff/handback timestamps for those same rows move
# forward a few seconds to reflect the freshly generated report.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 13, 14)

# zh-cn sheet: mark handoff priority and bump the "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-12 16:24:42"
}

# de-de sheet: mark handoff priority and bump the "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-12 16:24:51"
}

# Overview sheet: bump the "Latest HO Xliff Generate Date" column to match
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-12 16:24:51"
}
